# Rename the three logo pictures (wp:docPr / pic:cNvPr "name" attribute)
# so the display names match what they should be:
#   footer (default)    - Pearson logo : image2.png -> image1.png
#   footer (first page)  - Pearson logo : image2.png -> image1.png
#   header (first page)  - BTec logo    : image1.jpg -> image2.jpg
#
# InlineShape objects don't expose a settable "Name" in this object
# model, so each picture is temporarily converted to a floating Shape
# (which does support .Name), renamed, then converted back to an
# inline shape - this round-trip leaves every other property (size,
# id, wrapping, etc.) untouched.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlinePicture($range, $newName) {
    $inlineShapes = $range.InlineShapes
    if ($inlineShapes.Count -ge 1) {
        $pic = $inlineShapes.Item(1)
        $floating = $pic.ConvertToShape()
        $floating.Name = $newName
        $floating.ConvertToInlineShape() | Out-Null
    }
}

# Default (primary) footer - Pearson logo
Rename-InlinePicture $sec.Footers(1).Range "image1.png"

# First-page footer - Pearson logo
Rename-InlinePicture $sec.Footers(2).Range "image1.png"

# First-page header - BTec logo
Rename-InlinePicture $sec.Headers(2).Range "image2.jpg"
